$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 120, shifting existing rows 120-128 down to 121-129.
$ws.Rows("120:120").Insert()

# Populate the newly inserted row 120 with the new weekly record.
$ws.Range("A120").Value = 5
$ws.Range("B120").Value = "Macroferia Regional de Talca"
$ws.Range("C120").Value = "Maule"
$ws.Range("D120").Value = 45212
$ws.Range("E120").Value = 7
$ws.Range("F120").Value = 100112026
$ws.Range("G120").Value = "Haba"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 300
$ws.Range("K120").Value = 9000
$ws.Range("L120").Value = 9000
$ws.Range("M120").Value = 9000
$ws.Range("N120").Value = "`$/saco 25 kilos"
$ws.Range("O120").Value = "Región de O'Higgins"
$ws.Range("P120").Value = 360
$ws.Range("Q120").Value = 25
$ws.Range("R120").Value = "Hortaliza"
